# Weekly refresh of the "Fruta, Terminal Hortofrutícola Agro Chillán - Pera"
# sheet: the existing price rows (126-138) get this week's date/volume/price
# figures, and three additional observations (139-141) are appended that
# carry forward the previous week's figures that used to live at the tail
# of the range.
#
# NOTE: this runtime's `.Value` property getter is unreliable (it can
# surface the underlying COM property descriptor instead of the cell's
# content), so every read/write below goes through `.Value2`.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout (1-indexed): A Mercado ID, B Mercado, C Región, D Fecha,
# E Codreg, F Tipo, G Producto ID, H Producto, I Categoría ID, J Categoría,
# K Variedad, L Calidad, M Volumen, N Precio mínimo, O Precio máximo,
# P Precio promedio ponderado, Q Unidad, R Origen, S Precio $/Kg, T Kg/unidad.
$colD = 4
$colL = 12
$colM = 13
$colN = 14
$colO = 15
$colP = 16
$colS = 19

# Updated values for the existing rows 126-138.
$updates = @(
    @{ Row=126; D="2021-10-05"; L="Especial"; M=30;  N=11000; O=11000; P=11000; S=688 },
    @{ Row=127; D="2021-10-05"; L="Primera";  M=120; N=9500;  O=10000; P=9750;  S=609 },
    @{ Row=128; D="2021-10-05"; L="Segunda";  M=30;  N=9000;  O=9000;  P=9000;  S=562 },
    @{ Row=129; D="2021-08-30"; L="Especial"; M=120; N=9500;  O=10000; P=9750;  S=609 },
    @{ Row=130; D="2021-08-30"; L="Primera";  M=60;  N=8500;  O=9000;  P=8750;  S=547 },
    @{ Row=131; D="2021-03-31"; L="Primera";  M=120; N=9500;  O=10000; P=9750;  S=609 },
    @{ Row=132; D="2021-06-04"; L="Especial"; M=60;  N=11000; O=11000; P=11000; S=688 },
    @{ Row=133; D="2021-06-04"; L="Primera";  M=120; N=9000;  O=10000; P=9500;  S=594 },
    @{ Row=134; D="2021-04-20"; L="Primera";  M=120; N=10000; O=11000; P=10500; S=656 },
    @{ Row=135; D="2021-09-22"; L="Especial"; M=60;  N=11000; O=11000; P=11000; S=688 },
    @{ Row=136; D="2021-09-22"; L="Primera";  M=60;  N=9000;  O=10000; P=9500;  S=594 },
    @{ Row=137; D="2021-06-10"; L="Primera";  M=120; N=9000;  O=10000; P=9500;  S=594 },
    @{ Row=138; D="2021-06-10"; L="Segunda";  M=60;  N=8000;  O=8000;  P=8000;  S=500 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, $colD).Value2 = (Get-Date -Date $u.D)
    $ws.Cells.Item($r, $colL).Value2 = $u.L
    $ws.Cells.Item($r, $colM).Value2 = $u.M
    $ws.Cells.Item($r, $colN).Value2 = $u.N
    $ws.Cells.Item($r, $colO).Value2 = $u.O
    $ws.Cells.Item($r, $colP).Value2 = $u.P
    $ws.Cells.Item($r, $colS).Value2 = $u.S
}

# Three new rows (139-141) appended at the bottom, replicating the constant
# columns from row 138 and carrying the remaining observations forward.
$newRows = @(
    @{ Row=139; D="2021-08-17"; L="Especial"; M=80;  N=9500; O=10000; P=9750; S=609 },
    @{ Row=140; D="2021-08-17"; L="Primera";  M=60;  N=8500; O=9000;  P=8750; S=547 },
    @{ Row=141; D="2021-05-07"; L="Primera";  M=120; N=9000; O=10000; P=9500; S=594 }
)

$templateRow = 138
for ($col = 1; $col -le 20; $col++) {
    if ($col -eq $colD -or $col -eq $colL -or $col -eq $colM -or $col -eq $colN -or $col -eq $colO -or $col -eq $colP -or $col -eq $colS) {
        continue
    }
    $val = $ws.Cells.Item($templateRow, $col).Value2
    foreach ($n in $newRows) {
        $ws.Cells.Item($n.Row, $col).Value2 = $val
    }
}

foreach ($n in $newRows) {
    $r = $n.Row
    $ws.Cells.Item($r, $colD).Value2 = (Get-Date -Date $n.D)
    $ws.Cells.Item($r, $colL).Value2 = $n.L
    $ws.Cells.Item($r, $colM).Value2 = $n.M
    $ws.Cells.Item($r, $colN).Value2 = $n.N
    $ws.Cells.Item($r, $colO).Value2 = $n.O
    $ws.Cells.Item($r, $colP).Value2 = $n.P
    $ws.Cells.Item($r, $colS).Value2 = $n.S
}

# Keep the sheet's used-range dimension consistent with the new row count.
$ws.Cells.Item(141, 20).Value2 = $ws.Cells.Item(141, 20).Value2
